# Footer correction: "Sofi Wesson" <tab> <tab> "18/08/2022"
#                  -> "Sofi Wesson" <tab> <tab> "26/08/2022"
#
# In the authored edit, only the day-of-month digits were selected and
# retyped, so Word split the single run that held "18/08/2022" into two
# runs - "26" and "/08/2022" - rather than folding the new text back
# into one run. Everything else about the paragraph (the "Sofi Wesson"
# run, the two tab runs, and every attribute on the paragraph/runs) is
# left exactly as it was.

$d = $word.ActiveDocument

$footer = $d.Sections.Item(1).Footers.Item(1)
$footerRange = $footer.Range

$oldDate = "18/08/2022"
$newFirst = "26"
$newSecond = "/08/2022"

$footerText = $footerRange.Text
if ($footerText.IndexOf($oldDate) -lt 0) {
    throw "Could not find '$oldDate' in the footer (footer text was: '$footerText')"
}

# Pull the paragraph's current OOXML so the rewrite keeps every existing
# attribute (w14:paraId, w14:textId, the paragraph's w:rsidR/rsidRDefault,
# pPr, ...) untouched; we only touch the single run that holds the date.
$openXml = $footerRange.WordOpenXML
if (-not ($openXml -match '(?s)(<w:p\b.*?</w:p>)')) {
    throw "Could not locate the footer paragraph in WordOpenXML"
}
$paraXml = $Matches[1]

$oldRunXml = "<w:r><w:t>$oldDate</w:t></w:r>"
if ($paraXml.IndexOf($oldRunXml) -lt 0) {
    throw "Could not find the date run '$oldRunXml' in paragraph XML: $paraXml"
}
$newRunXml = "<w:r><w:t>$newFirst</w:t></w:r><w:r><w:t>$newSecond</w:t></w:r>"
$newParaXml = $paraXml.Replace($oldRunXml, $newRunXml)

# WordOpenXML normalises same-formatting runs (e.g. the two plain <w:tab/>
# runs) and drops their w:rsidR when it serialises them back out. Those
# two tab runs originally each carried w:rsidR="00853E7F"; restore it so
# the round trip doesn't silently erase that (unrelated, unchanged)
# attribute.
$newParaXml = $newParaXml -replace '<w:r><w:tab/></w:r>', '<w:r w:rsidR="00853E7F"><w:tab/></w:r>'

# The extracted fragment needs its own xmlns:w declaration to stand alone
# as input to InsertXML.
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'
$newParaXml = $newParaXml -replace '^<w:p ', "<w:p $wNs "

# InsertXML only behaves predictably when the target Range is the entire
# paragraph (a sub-range clears the whole paragraph instead of patching
# it), so replace the whole footer paragraph range in one call.
$footerRange.InsertXML($newParaXml) | Out-Null

$finalText = $footer.Range.Text
if ($finalText.IndexOf("26/08/2022") -lt 0) {
    throw "Footer did not update as expected, got: '$finalText'"
}

Write-Output "Footer now reads: $finalText"
